# B0CHSGZPNY_po_data.xlsx — add PO Forecast sheet + rename headers
$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" header on the existing sheets -----
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet as the last tab --------------------
# Duplicate "Weekly Quantity" (carries over sheetPr/pageMargins/formats) and
# place the copy after the last existing sheet, then rename + clean it up.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsWeekly.Copy($null, $lastSheet)
$wsForecast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Name = "PO Forecast"

# Wipe the template's old data/values, keep column A's date formatting and
# extend the header format / date format across the new C:D columns and
# down through row 29 (date-formatted "A" column, s="2").
$wsForecast.Cells.ClearContents()
$wsForecast.Range("B1").Copy()
$wsForecast.Range("C1:D1").PasteSpecial(-4122)
$wsForecast.Range("A2").Copy()
$wsForecast.Range("A2:A29").PasteSpecial(-4122)

# --- 3. Header row ----------------------------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 4. Forecast data rows 2-29 ---------------------------------------------
$poData = New-Object 'object[,]' 28,4

$poData[0,0] = 45193.99999999999; $poData[0,1] = 60; $poData[0,2] = 5.558279708567823; $poData[0,3] = 112.0076268545824
$poData[1,0] = 45214.99999999999; $poData[1,1] = 57; $poData[1,2] = 2.596379653283505; $poData[1,3] = 109.0268725645244
$poData[2,0] = 45221.99999999999; $poData[2,1] = 55; $poData[2,2] = 3.207175380011305; $poData[2,3] = 105.0965991999008
$poData[3,0] = 45228.99999999999; $poData[3,1] = 54; $poData[3,2] = 1.897186991850533; $poData[3,3] = 106.8651288006023
$poData[4,0] = 45242.99999999999; $poData[4,1] = 52; $poData[4,2] = -1.671874767809168; $poData[4,3] = 102.9479751030712
$poData[5,0] = 45249.99999999999; $poData[5,1] = 50; $poData[5,2] = -2.134522490339942; $poData[5,3] = 99.39734584205654
$poData[6,0] = 45256.99999999999; $poData[6,1] = 49; $poData[6,2] = -7.093632256712498; $poData[6,3] = 99.85199125738538
$poData[7,0] = 45263.99999999999; $poData[7,1] = 48; $poData[7,2] = -6.362434797018753; $poData[7,3] = 101.4309804226576
$poData[8,0] = 45270.99999999999; $poData[8,1] = 47; $poData[8,2] = -5.434581302484518; $poData[8,3] = 101.8323660400138
$poData[9,0] = 45277.99999999999; $poData[9,1] = 45; $poData[9,2] = -4.306309085361154; $poData[9,3] = 98.58005129001768
$poData[10,0] = 45298.99999999999; $poData[10,1] = 42; $poData[10,2] = -9.891658587326662; $poData[10,3] = 90.89208041007365
$poData[11,0] = 45305.99999999999; $poData[11,1] = 40; $poData[11,2] = -13.54169095771115; $poData[11,3] = 88.75878274626047
$poData[12,0] = 45312.99999999999; $poData[12,1] = 39; $poData[12,2] = -15.42366885975015; $poData[12,3] = 91.07167847282194
$poData[13,0] = 45319.99999999999; $poData[13,1] = 38; $poData[13,2] = -15.47202824889044; $poData[13,3] = 88.12006151299231
$poData[14,0] = 45326.99999999999; $poData[14,1] = 36; $poData[14,2] = -17.58150351437851; $poData[14,3] = 88.5082124314991
$poData[15,0] = 45333.99999999999; $poData[15,1] = 35; $poData[15,2] = -15.40606922167427; $poData[15,3] = 90.53304301192296
$poData[16,0] = 45340.99999999999; $poData[16,1] = 34; $poData[16,2] = -16.3900433209003; $poData[16,3] = 85.16234347367121
$poData[17,0] = 45347.99999999999; $poData[17,1] = 33; $poData[17,2] = -19.38518063669054; $poData[17,3] = 87.35430819617179
$poData[18,0] = 45396.99999999999; $poData[18,1] = 24; $poData[18,2] = -32.37170632712959; $poData[18,3] = 75.77377585886271
$poData[19,0] = 45403.99999999999; $poData[19,1] = 23; $poData[19,2] = -27.63602673789943; $poData[19,3] = 73.49182855279184
$poData[20,0] = 45410.99999999999; $poData[20,1] = 21; $poData[20,2] = -29.85782300734045; $poData[20,3] = 67.85004559853311
$poData[21,0] = 45417.99999999999; $poData[21,1] = 20; $poData[21,2] = -30.97711386539762; $poData[21,3] = 67.48485735007301
$poData[22,0] = 45424.99999999999; $poData[22,1] = 19; $poData[22,2] = -35.78715961290863; $poData[22,3] = 71.60259577971645
$poData[23,0] = 45431.99999999999; $poData[23,1] = 18; $poData[23,2] = -34.05409393142838; $poData[23,3] = 67.88971457437428
$poData[24,0] = 45438.99999999999; $poData[24,1] = 16; $poData[24,2] = -34.8639759672021; $poData[24,3] = 67.8006844241072
$poData[25,0] = 45445.99999999999; $poData[25,1] = 15; $poData[25,2] = -39.26663849780242; $poData[25,3] = 65.76623733133474
$poData[26,0] = 45452.99999999999; $poData[26,1] = 14; $poData[26,2] = -40.69694557174985; $poData[26,3] = 59.79336848891928
$poData[27,0] = 45459.99999999999; $poData[27,1] = 13; $poData[27,2] = -43.02349905599796; $poData[27,3] = 61.18619836880307

$wsForecast.Range("A2:D29").Value = $poData

# Restore the original active sheet/selection (first tab, cell A1)
$wsWeekly.Select() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
